$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 25-37 (cell-level changes) ---
# Row 25
$ws.Range("D25").Value = 44435

# Row 26
$ws.Range("D26").Value = 44431
$ws.Range("J26").Value = 120
$ws.Range("K26").Value = 28000
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = 29000
$ws.Range("P26").Value = 1933

# Row 27
$ws.Range("D27").Value = 44330
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 28000
$ws.Range("L27").Value = 30000
$ws.Range("M27").Value = 29000
$ws.Range("P27").Value = 1933

# Row 28
$ws.Range("D28").Value = 44358
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 21000
$ws.Range("L28").Value = 22000
$ws.Range("M28").Value = 21500
$ws.Range("P28").Value = 1433

# Row 29
$ws.Range("D29").Value = 44389
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 21000
$ws.Range("L29").Value = 22000
$ws.Range("M29").Value = 21500
$ws.Range("P29").Value = 1433

# Row 30
$ws.Range("D30").Value = 44399
$ws.Range("I30").Value = "Primera"
$ws.Range("K30").Value = 20000
$ws.Range("L30").Value = 21000
$ws.Range("M30").Value = 20500
$ws.Range("P30").Value = 1367

# Row 31
$ws.Range("D31").Value = 44305
$ws.Range("K31").Value = 13000
$ws.Range("L31").Value = 14000
$ws.Range("M31").Value = 13500
$ws.Range("P31").Value = 900

# Row 32
$ws.Range("D32").Value = 44333
$ws.Range("H32").Value = "Inferno"
$ws.Range("I32").Value = "Segunda"

# Row 33
$ws.Range("D33").Value = 44309
$ws.Range("H33").Value = "Inferno"
$ws.Range("I33").Value = "Primera"
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 11500
$ws.Range("P33").Value = 767

# Row 34
$ws.Range("D34").Value = 44344
$ws.Range("H34").Value = "Cristal"
$ws.Range("K34").Value = 24000
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = 24500
$ws.Range("P34").Value = 1633

# Row 35
$ws.Range("D35").Value = 44316
$ws.Range("H35").Value = "Cristal"
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 9500
$ws.Range("P35").Value = 633

# Row 36
$ws.Range("D36").Value = 44232
$ws.Range("K36").Value = 17000
$ws.Range("L36").Value = 18000
$ws.Range("M36").Value = 17500
$ws.Range("P36").Value = 1167

# Row 37
$ws.Range("D37").Value = 44232
$ws.Range("I37").Value = "Segunda"
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 16000
$ws.Range("M37").Value = 15500
$ws.Range("P37").Value = 1033

# --- Add new rows 38-39 ---
# Row 38
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C38").Value = "Arica y Parinacota"
$ws.Range("D38").Value = 44270
$ws.Range("E38").Value = 15
$ws.Range("F38").Value = 100112021
$ws.Range("G38").Value = "Ají"
$ws.Range("H38").Value = "Inferno"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 120
$ws.Range("K38").Value = 12000
$ws.Range("L38").Value = 13000
$ws.Range("M38").Value = 12500
$ws.Range("N38").Value = "$/caja 15 kilos"
$ws.Range("O38").Value = "Región de Arica y Parinacota"
$ws.Range("P38").Value = 833
$ws.Range("Q38").Value = 15
$ws.Range("R38").Value = "Hortaliza"
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 39
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C39").Value = "Arica y Parinacota"
$ws.Range("D39").Value = 44284
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = 100112021
$ws.Range("G39").Value = "Ají"
$ws.Range("H39").Value = "Inferno"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 120
$ws.Range("K39").Value = 11000
$ws.Range("L39").Value = 12000
$ws.Range("M39").Value = 11500
$ws.Range("N39").Value = "$/caja 15 kilos"
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 767
$ws.Range("Q39").Value = 15
$ws.Range("R39").Value = "Hortaliza"
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"

